$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update state / name / gender values in the data rows to match the
# refactored importer test fixture (order matches the shared-string
# append order seen in the target workbook).
$ws.Range("D3").Value  = "Hamburg"
$ws.Range("D4").Value  = "centre"
$ws.Range("D8").Value  = "Buckinghamshire"
$ws.Range("D9").Value  = "Liverpool"
$ws.Range("A12").Value = "Hyphenated"
$ws.Range("A7").Value  = "Michael"
$ws.Range("B7").Value  = "Baldwin"
$ws.Range("C4").Value  = "Male"
$ws.Range("C9").Value  = "FEM"
$ws.Range("C13").Value = "male"

# Move the active selection from H17 to C1.
$ws.Range("C1").Select()
